# "Scanned. Before Repeater Test"
# Appends one new data row (row 10) to the sheet, mirroring the previous
# row's layout, and fills in a new score column (X/Y) for the existing
# last row (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 10 by copying row 9 (values + number formats/styles) so the
# new date cell (A10) and the percentage cells (S10/T10) inherit the same
# formatting as the rest of the table instead of picking up Excel's
# General default.
$ws.Range("A9:W9").Copy($ws.Range("A10"))

# New X/Y score columns for the previously-last row.
$ws.Range("X9").Value = 0.47000100000001055
$ws.Range("Y9").Value = "Up"

# Newly appended row of scan data.
$ws.Range("A10").Value = 42649.886817129627
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Neutral"
$ws.Range("D10").Value = 22
$ws.Range("E10").Value = 19184
$ws.Range("F10").Value = 2333
$ws.Range("G10").Value = 54
$ws.Range("H10").Value = 44
$ws.Range("I10").Value = 84
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 10901
$ws.Range("L10").Value = 239
$ws.Range("M10").Value = 194
$ws.Range("N10").Value = 99
$ws.Range("O10").Value = 18
$ws.Range("P10").Value = "Bag"
$ws.Range("Q10").Value = 35.958706302092025
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = -0.0019
$ws.Range("T10").Value = -0.023
$ws.Range("U10").Value = 14.71
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = 0
